# Trade #13 closed at 2026-02-17 15:17:33 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.75   # Current Capital
$wsSummary.Range("B4").Value = -0.25     # Total P&L $
$wsSummary.Range("B5").Value = -0.38     # Total P&L %
$wsSummary.Range("B6").Value = 13        # Total Trades
$wsSummary.Range("B8").Value = 6         # Losing Trades
$wsSummary.Range("B9").Value = 23.08     # Win Rate %

# --- Strategy Status sheet (MarketMaking row) ---
$wsStrategy = $wb.Worksheets.Item("Strategy Status")
$wsStrategy.Range("C4").Value = 99.75    # Capital
$wsStrategy.Range("D4").Value = 13       # Trades
$wsStrategy.Range("E4").Value = -0.25    # P&L $
$wsStrategy.Range("F4").Value = -0.25    # P&L %
$wsStrategy.Range("G4").Value = 23.08    # Win Rate %

# --- Append new trade row (#13 -> row 14) to "All Trades" and "MarketMaking" sheets ---
$sheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A14").Value = 13

    # The date column is stored as plain text, not an Excel date serial,
    # so force text format before assignment to avoid auto-conversion to
    # a date number, then restore the default "Normal" style so the cell
    # keeps the same (unstyled) look as the rest of the sheet.
    $ws.Range("B14").NumberFormat = "@"
    $ws.Range("B14").Value = "2026-02-17"
    $ws.Range("B14").Style = "Normal"
    $ws.Range("C14").Value = "15:17:25"

    $ws.Range("D14").Value = "MarketMaking"
    $ws.Range("E14").Value = "UP"
    $ws.Range("F14").Value = 0.9061940000000001
    $ws.Range("G14").Value = 0.786759
    $ws.Range("H14").Value = "CLOSED"
    $ws.Range("I14").Value = -13.1799
    $ws.Range("J14").Value = -0.12
    $ws.Range("K14").Value = 99.75
    $ws.Range("L14").Value = 0
    $ws.Range("M14").Value = 0
    $ws.Range("N14").Value = 0.6
    $ws.Range("O14").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P14").Value = "early_exit"
    $ws.Range("Q14").Value = 0.12
}

Write-Host "Applied trade #13 updates"
